$d = $word.ActiveDocument

function Replace-RangeXml([string]$searchText, [string]$innerXml) {
    $rngFind = $d.Content.Duplicate
    $found = $rngFind.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    # Re-seat the match as a fresh Range; Find-mutated ranges misbehave under InsertXML.
    $rng = $d.Range($rngFind.Start, $rngFind.End)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# 1) Split "J. Holwerda" into "J. " + proofErr-wrapped "Holwerda"
$inner1 = '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">J. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Holwerda</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Replace-RangeXml "J. Holwerda" $inner1

# 2) Restructure the "upload" / "surfspot" paragraph text
$searchText2 = ", upload functie en de nieuws pagina. De coach was erg te spreken over onze vorderingen en vond dat wij op de goede weg waren. Verder hebben we het er over gehad dat wij deze week bezig gaan met het LDAP systeem dit houdt in dat wij ons gaan inlezen in dit onderwerp en zo snel mogelijk een koppeling gaan maken met het bestaand systeem. We hadden echter nog een vraag over het login systeem op de site van surfspot hadden we gezien dat deze website een connectie heeft met "
$inner2 = '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>upload</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> functie en de nieuws pagina. De coach was erg te spreken over onze vorderingen en vond dat wij op de goede weg waren.</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> Gevraagd naar hoe het Technisch ontwerp eruit moet komen te zien.</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> Verder hebben we het er over gehad dat wij deze week bezig gaan met het LDAP systeem dit houdt in dat wij ons gaan inlezen in dit onderwerp en zo snel mogelijk een koppeling gaan maken met het bestaand systeem. We hadden echter nog een vraag over het login systeem op de site van surfspot hadden we gezien dat deze website een connectie heeft met </w:t></w:r>'
Replace-RangeXml $searchText2 $inner2
